# edit.ps1 - apply the "Celestial Symphony" -> "Allure of Chemistry" rewrite
# described by the target diff, via Word COM-interop calls.

$d = $word.ActiveDocument

# Replace $old with $new, scoped to paragraph number $paraNum (1-based) to
# avoid ambiguous matches elsewhere in the document.
function Replace-In-Para([int]$paraNum, [string]$old, [string]$new) {
    $p = $d.Paragraphs($paraNum)
    $rng = $p.Range.Duplicate
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARN: replace failed in para $paraNum for: $old"
    }
    return $ok
}

# Forces the engine to split the run that currently spans [start, start+len)
# away from its neighbours, without altering the visible formatting, by
# toggling a character property on and back off.
function Split-Run([int]$start, [int]$len) {
    $r = $d.Range($start, $start + $len)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# Inserts $text right after position $pos, then splits it into its own run.
# Returns the position right after the inserted text.
function Insert-Run([int]$pos, [string]$text) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($text)
    Split-Run $pos $text.Length
    return $pos + $text.Length
}

# Finds $needle inside paragraph $paraNum and returns its End offset.
function Find-End-In-Para([int]$paraNum, [string]$needle) {
    $p = $d.Paragraphs($paraNum)
    $rng = $p.Range.Duplicate
    $rng.Find.Execute($needle) | Out-Null
    return $rng.End
}

# Finds $needle inside paragraph $paraNum and returns its Start offset.
function Find-Start-In-Para([int]$paraNum, [string]$needle) {
    $p = $d.Paragraphs($paraNum)
    $rng = $p.Range.Duplicate
    $rng.Find.Execute($needle) | Out-Null
    return $rng.Start
}

# ---------------------------------------------------------------------------
# Title (paragraph 1)
# ---------------------------------------------------------------------------
Replace-In-Para 1 "Celestial Symphony: Unraveling the Harmony of Stellar Sounds" "The Allure of Chemistry: Unraveling the Molecular Symphony of Life"

# ---------------------------------------------------------------------------
# Author name (paragraph 2)
# ---------------------------------------------------------------------------
Replace-In-Para 2 " Adriana Anderson" " Julia Carter"

# ---------------------------------------------------------------------------
# Email line (paragraph 3) -> collapses to a single run containing "at"
# ---------------------------------------------------------------------------
Replace-In-Para 3 "andriana.astro@heavenlyharmonic.com" "at"

# ---------------------------------------------------------------------------
# Body paragraph (paragraph 5)
# ---------------------------------------------------------------------------
Replace-In-Para 5 "The universe, an infinite expanse of cosmic wonders, holds secrets that have captivated humankind for centuries" "In the vast tapestry of scientific disciplines, chemistry stands as a beacon of wonder, revealing the intricate dance of atoms and molecules that orchestrates the symphony of life"

Replace-In-Para 5 " Among these mysteries lies the enigmatic realm of stellar sounds, a symphony of cosmic vibrations that has long intrigued scientists and music enthusiasts alike" " Chemistry is the study of matter, both living and nonliving, and their physical and chemical properties"

# This run splits into three: itself (shortened), a new ".", and a new sentence.
Replace-In-Para 5 " As we delve into the celestial realm, we embark on a journey to unravel the harmony of stellar sounds, exploring the mechanisms that produce these cosmic melodies and the profound insights they offer into the workings of our universe" " It explores the interactions between substances, uncovering the secrets of their composition, structure, and reactivity"
$pos = Find-End-In-Para 5 " It explores the interactions between substances, uncovering the secrets of their composition, structure, and reactivity"
$pos = Insert-Run $pos "."
$pos = Insert-Run $pos " As we delve into the fascinating world of chemistry, we embark on a journey that unveils the fundamental principles governing the natural world, promising a deeper understanding of the intricate processes that underpin life"

Replace-In-Para 5 "The celestial symphony begins with the mesmerizing song of stars, born amidst the fiery dance of interstellar gas and dust" "From the everyday phenomena we witness to the complex reactions occurring within our bodies, chemistry holds the key to comprehending the remarkable diversity and unity observed in the universe"

Replace-In-Para 5 " As these incandescent orbs fuse elements in their cores, they emit pulsations that ripple through space as vibrations" " With its focus on matter and its transformations, chemistry offers a lens through which we can decipher the enigmas of the molecular world"

Replace-In-Para 5 " These stellar oscillations, driven by complex interactions between gravity, pressure, and temperature, produce distinct patterns of sound, each star possessing its unique celestial tune" " It unveils the secrets of chemical reactions, revealing how substances combine, rearrange, and decompose, opening up new avenues for innovation in medicine, energy, and materials science"

Replace-In-Para 5 "Venturing beyond the solitary symphonies of individual stars, we encounter celestial choirs where entire star clusters unite to create cosmic harmonies" "Chemistry is an experimental science"

Replace-In-Para 5 " Gravitational interactions between stars within these clusters give rise to collective oscillations, resulting in synchronized pulsations and rhythmic modulations that resonate across vast distances" " By skillfully conducting experiments, chemists meticulously investigate the properties of substances, unraveling their innermost secrets"

# This run splits into three as well.
Replace-In-Para 5 " The combined sounds of these stellar ensembles weave intricate tapestries of sound, adding depth and complexity to the celestial symphony" " Through careful observation, analysis, and interpretation of data, they uncover the mechanisms underlying chemical reactions, shedding light on the intricate choreography of atoms and molecules"
$pos = Find-End-In-Para 5 " Through careful observation, analysis, and interpretation of data, they uncover the mechanisms underlying chemical reactions, shedding light on the intricate choreography of atoms and molecules"
$pos = Insert-Run $pos "."
$pos = Insert-Run $pos " The profound insights gained from these experiments have shaped our understanding of the universe, laying the foundation for groundbreaking advancements that have transformed society"

# ---------------------------------------------------------------------------
# Summary paragraph (paragraph 7)
# ---------------------------------------------------------------------------
Replace-In-Para 7 "The exploration of stellar sounds unveils a hidden realm of cosmic beauty and scientific intrigue" "Chemistry is the study of matter and its properties, revealing the intricate molecular symphony of life"

# The run containing the old "From the solitary tunes..." text becomes two
# runs; the second one carries the relocated <w:lastRenderedPageBreak/>.
Replace-In-Para 7 " From the solitary tunes of individual stars to the celestial harmonies of star clusters, the universe reverberates with a symphony of vibrations" " It explores the composition, structure, and reactivity of substances, "
$pos = Find-End-In-Para 7 " It explores the composition, structure, and reactivity of substances, "
$pos = Insert-Run $pos "unveiling the secrets of chemical reactions and unlocking the enigmas of the molecular world"

Replace-In-Para 7 " These sounds, born from the fundamental forces that govern the cosmos, offer unique insights into the enigmatic lives of stars, the dynamics of stellar clusters, and the vast symphony of the universe" " With its focus on experimentation and meticulous data analysis, chemistry provides a deeper understanding of the fundamental principles governing the natural world and paves the way for groundbreaking advancements in various fields"

# The old " As " run and the following lastRenderedPageBreak run merge into
# a single new run (the page break marker moves earlier in the paragraph,
# so it no longer belongs here).
$startPos = Find-Start-In-Para 7 " As "
$endPos = Find-End-In-Para 7 "we continue to unravel the mysteries of stellar sounds, we deepen our understanding of the harmonious interplay between cosmic phenomena and the fundamental principles that underpin the fabric of our universe"
$oldSpan = $d.Range($startPos, $endPos)
$oldSpan.Text = " Through chemistry, we gain insights that not only unravel the mysteries of the universe but also hold the potential to transform society and improve human lives"

# ---------------------------------------------------------------------------
# New trailing empty paragraph before the final section break.
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "done"
